$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "32.994.55"
$ws.Range("E2").Value = "  +10.22%  "
$ws.Range("D3").Value = "1.762.39"
$ws.Range("E3").Value = "  +5.99%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.543"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("E10").Value = "  +5.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0666"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0922"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "2.014.42"
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("D14").Value = "1.787.58"
$ws.Range("E14").Value = "  +7.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.41%  "
$ws.Range("D18").Value = "32.904.63"
$ws.Range("E18").Value = "  +9.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "260.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.29%  "
$ws.Range("D21").Value = "0.0₃0747"
$ws.Range("E21").Value = "  +4.95%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.12%  "
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.51%  "
$ws.Range("E28").Value = "  +4.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0515"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.99%  "
$ws.Range("D35").Value = "1.554.07"
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "85.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.37%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.625"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0186"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.883"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("E44").Value = "  +6.00%  "
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.77%  "
$ws.Range("D48").Value = "1.915.45"
$ws.Range("E48").Value = "  +5.88%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.95%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +21.98%  "
